# Prefab change for strings
# Delete the bug-report row about the lever/hook wiggle issue from the
# "Bug Tracking" sheet (row 13), shifting all following rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Tracking")
$ws.Activate()

# Delete the entire row 13 (shifts rows 14+ up to 13+, updates
# dimension, conditional formatting and data validation ranges
# automatically, and Excel will drop now-unused shared strings on save).
$ws.Rows.Item(13).Delete()

# Reflect the new selection state left in the saved file.
$ws.Range("B11").Select()
